# BookingConfirmation.xlsx - "excel list and qty fix"
#
# The "Booking Confirmation" sheet has a set of D-column cells that hold a
# stale/broken `=#REF!` formula (left over from a deleted reference sheet).
# Clear those out so the cells go back to being plain (empty) cells while
# keeping their existing number formatting/style untouched, then leave the
# final selection on D57 (the last of the cleaned-up cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Booking Confirmation")

$brokenRefCells = @(
    "D15",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D27",
    "D29",
    "D33",
    "D49",
    "D57"
)

foreach ($cellRef in $brokenRefCells) {
    $ws.Range($cellRef).ClearContents()
}

# Match the final on-sheet selection left behind by the edit.
$ws.Range("D57").Select()
